# Applies the "Updated cryptos list ... with GitHub Actions" refresh:
# per-row Price/Volume(1h) updates, plus a rank swap between the
# Bittensor and dogwifhat rows (39/40).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin / Link / Volume(1h) text cells: plain assignment is safe because
# none of these strings is parsed by Excel as a pure number (percentages
# carry padding spaces, the dotted prices have 2+ separators, etc).
$ws.Range('D2').Value = '64.307.89'
$ws.Range('D3').Value = '3.162.68'
$ws.Range('E3').Value = '  -0.71%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  +1.75%  '
$ws.Range('E6').Value = '  -2.81%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').Value = '3.157.26'
$ws.Range('E8').Value = '  -0.91%  '
$ws.Range('E9').Value = '  -0.58%  '
$ws.Range('E10').Value = '  -0.63%  '
$ws.Range('E11').Value = '  -1.85%  '
$ws.Range('E12').Value = '  -0.60%  '
$ws.Range('E13').Value = '  +1.53%  '
$ws.Range('E14').Value = '  -3.41%  '
$ws.Range('D15').Value = '3.684.11'
$ws.Range('E15').Value = '  -0.65%  '
$ws.Range('E16').Value = '  +3.10%  '
$ws.Range('D17').Value = '64.270.99'
$ws.Range('E17').Value = '  -0.85%  '
$ws.Range('D18').Value = '3.164.74'
$ws.Range('E18').Value = '  -0.59%  '
$ws.Range('E19').Value = '  -1.68%  '
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('E22').Value = '  +0.60%  '
$ws.Range('E24').Value = '  -0.76%  '
$ws.Range('E25').Value = '  -1.26%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('E27').Value = '  -3.10%  '
$ws.Range('E28').Value = '  -0.38%  '
$ws.Range('E29').Value = '  +3.42%  '
$ws.Range('E30').Value = '  -5.28%  '
$ws.Range('E31').Value = '  -6.27%  '
$ws.Range('E32').Value = '  +0.39%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('E34').Value = '  -1.81%  '
$ws.Range('E35').Value = '  +1.18%  '
$ws.Range('D36').Value = '0.0₃0801'
$ws.Range('E36').Value = '  +8.62%  '
$ws.Range('E37').Value = '  -1.70%  '
$ws.Range('E38').Value = '  -2.50%  '
$ws.Range('B39').Value = 'dogwifhat'
$ws.Range('C39').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('E39').Value = '  -2.23%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('E40').Value = '  +1.32%  '
$ws.Range('E41').Value = '  -0.55%  '
$ws.Range('E42').Value = '  -4.22%  '
$ws.Range('E43').Value = '  -1.55%  '
$ws.Range('D44').Value = '2.875.54'
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('E45').Value = '  -4.37%  '
$ws.Range('E46').Value = '  -1.84%  '
$ws.Range('E47').Value = '  +3.96%  '
$ws.Range('E48').Value = '  -2.70%  '
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('E50').Value = '  -1.75%  '
$ws.Range('E51').Value = '  -1.05%  '

# Price cells that look like plain decimal numbers (e.g. "610.87",
# "0.720") must stay stored as text, matching the original inline-string
# cell type/format (General, no explicit style). Assigning .Value directly
# would have Excel auto-convert them to numbers, so instead we stage the
# text in an unused helper cell (A1, blank in the source sheet) that is
# explicitly formatted as Text, copy it, and Paste Special "Values only"
# (xlPasteValues = -4163) into the target - this carries over just the text,
# leaving the target cell's own style untouched. The helper cell is fully
# reset afterwards with Clear() so it is left exactly as it started (blank,
# default style).
$helper = $ws.Range('A1')
$helper.NumberFormat = '@'

$helper.Value = '610.87'
$helper.Copy()
$ws.Range('D5').PasteSpecial(-4163)
$helper.Value = '147.83'
$helper.Copy()
$ws.Range('D6').PasteSpecial(-4163)
$helper.Value = '0.0000262'
$helper.Copy()
$ws.Range('D13').PasteSpecial(-4163)
$helper.Value = '35.65'
$helper.Copy()
$ws.Range('D14').PasteSpecial(-4163)
$helper.Value = '479.45'
$helper.Copy()
$ws.Range('D20').PasteSpecial(-4163)
$helper.Value = '14.76'
$helper.Copy()
$ws.Range('D21').PasteSpecial(-4163)
$helper.Value = '0.720'
$helper.Copy()
$ws.Range('D22').PasteSpecial(-4163)
$helper.Value = '8.05'
$helper.Copy()
$ws.Range('D23').PasteSpecial(-4163)
$helper.Value = '13.81'
$helper.Copy()
$ws.Range('D24').PasteSpecial(-4163)
$helper.Value = '83.80'
$helper.Copy()
$ws.Range('D25').PasteSpecial(-4163)
$helper.Value = '2.84'
$helper.Copy()
$ws.Range('D27').PasteSpecial(-4163)
$helper.Value = '8.62'
$helper.Copy()
$ws.Range('D28').PasteSpecial(-4163)
$helper.Value = '7.20'
$helper.Copy()
$ws.Range('D29').PasteSpecial(-4163)
$helper.Value = '2.74'
$helper.Copy()
$ws.Range('D32').PasteSpecial(-4163)
$helper.Value = '26.43'
$helper.Copy()
$ws.Range('D34').PasteSpecial(-4163)
$helper.Value = '6.03'
$helper.Copy()
$ws.Range('D37').PasteSpecial(-4163)
$helper.Value = '53.17'
$helper.Copy()
$ws.Range('D38').PasteSpecial(-4163)
$helper.Value = '3.19'
$helper.Copy()
$ws.Range('D39').PasteSpecial(-4163)
$helper.Value = '464.49'
$helper.Copy()
$ws.Range('D40').PasteSpecial(-4163)
$helper.Value = '0.0401'
$helper.Copy()
$ws.Range('D41').PasteSpecial(-4163)
$helper.Value = '2.33'
$helper.Copy()
$ws.Range('D45').PasteSpecial(-4163)
$helper.Value = '0.271'
$helper.Copy()
$ws.Range('D46').PasteSpecial(-4163)
$helper.Value = '2.44'
$helper.Copy()
$ws.Range('D47').PasteSpecial(-4163)
$helper.Value = '0.999'
$helper.Copy()
$ws.Range('D49').PasteSpecial(-4163)
$helper.Value = '119.19'
$helper.Copy()
$ws.Range('D51').PasteSpecial(-4163)

$helper.Clear()
$excel.CutCopyMode = $false

